$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A16").Value = "wfwef"
$ws.Range("A17").Value = "wfwfwef"
$ws.Range("A18").Value = "wefwe"
$ws.Range("A19").Value = "wefwe"
$ws.Range("A20").Value = "fwef"
$ws.Range("A21").Value = "fwefwe"

$ws.Range("A22").Select()
